$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing date values (column D) for rows 126-135 following the
# described shift pattern.
$ws.Range("D126").Value = 44476
$ws.Range("D127").Value = 44476
$ws.Range("D128").Value = 44386
$ws.Range("D129").Value = 44386
$ws.Range("D130").Value = 44306
$ws.Range("D131").Value = 44306
$ws.Range("D132").Value = 44425
$ws.Range("D133").Value = 44425
$ws.Range("D134").Value = 44187
$ws.Range("D135").Value = 44187

# Append two new rows (136, 137) that replicate rows 134/135's full
# content but with the "new" date value (44250) that was displaced by the
# shift above.
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 44250
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112040
$ws.Range("G136").Value = "Cilantro"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 200
$ws.Range("K136").Value = 600
$ws.Range("L136").Value = 700
$ws.Range("M136").Value = 650
$ws.Range("N136").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O136").Value = "Región de Ñuble"
$ws.Range("P136").Value = 650
$ws.Range("Q136").Value = 1
$ws.Range("R136").Value = "Hortaliza"

$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44250
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112040
$ws.Range("G137").Value = "Cilantro"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 100
$ws.Range("K137").Value = 500
$ws.Range("L137").Value = 500
$ws.Range("M137").Value = 500
$ws.Range("N137").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O137").Value = "Región de Ñuble"
$ws.Range("P137").Value = 500
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"

# Match the style (date number format) that column D uses elsewhere in
# the sheet for the two newly appended rows.
$ws.Range("D136").NumberFormat = $ws.Range("D135").NumberFormat
$ws.Range("D137").NumberFormat = $ws.Range("D135").NumberFormat
